$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I10").Value = 'sd'
$ws.Range("J10").Value = 'Statement-non-opinion'
$ws.Range("I26").Value = 'sd'
$ws.Range("J26").Value = 'Statement-non-opinion'
$ws.Range("I31").Value = 'sd'
$ws.Range("J31").Value = 'Statement-non-opinion'
$ws.Range("I41").Value = 'sd'
$ws.Range("J41").Value = 'Statement-non-opinion'
$ws.Range("I52").Value = 'sv'
$ws.Range("J52").Value = 'Statement-opinion'
$ws.Range("I54").Value = 'aa'
$ws.Range("J54").Value = 'Agree/Accept'
$ws.Range("I70").Value = 'aa'
$ws.Range("J70").Value = 'Agree/Accept'
$ws.Range("I74").Value = 'aa'
$ws.Range("J74").Value = 'Agree/Accept'
$ws.Range("I88").Value = 'aa'
$ws.Range("J88").Value = 'Agree/Accept'
$ws.Range("I89").Value = 'aa'
$ws.Range("J89").Value = 'Agree/Accept'
$ws.Range("I120").Value = 'aa'
$ws.Range("J120").Value = 'Agree/Accept'
$ws.Range("I123").Value = 'aa'
$ws.Range("J123").Value = 'Agree/Accept'
$ws.Range("I134").Value = 'sd'
$ws.Range("J134").Value = 'Statement-non-opinion'
$ws.Range("I136").Value = 'aa'
$ws.Range("J136").Value = 'Agree/Accept'
$ws.Range("I137").Value = 'sd'
$ws.Range("J137").Value = 'Statement-non-opinion'
$ws.Range("I150").Value = 'sd'
$ws.Range("J150").Value = 'Statement-non-opinion'
$ws.Range("I154").Value = 'b'
$ws.Range("J154").Value = 'Acknowledge (Backchannel)'
$ws.Range("I170").Value = 'sv'
$ws.Range("J170").Value = 'Statement-opinion'
$ws.Range("I178").Value = 'sd'
$ws.Range("J178").Value = 'Statement-non-opinion'
$ws.Range("I194").Value = 'sd'
$ws.Range("J194").Value = 'Statement-non-opinion'
$ws.Range("I202").Value = 'sd'
$ws.Range("J202").Value = 'Statement-non-opinion'
$ws.Range("I209").Value = 'sd'
$ws.Range("J209").Value = 'Statement-non-opinion'
$ws.Range("I210").Value = 'b'
$ws.Range("J210").Value = 'Acknowledge (Backchannel)'
$ws.Range("I211").Value = 'sv'
$ws.Range("J211").Value = 'Statement-opinion'
$ws.Range("I214").Value = '%'
$ws.Range("J214").Value = 'Uninterpretable'
$ws.Range("I215").Value = '%'
$ws.Range("J215").Value = 'Uninterpretable'
$ws.Range("I217").Value = 'aa'
$ws.Range("J217").Value = 'Agree/Accept'
$ws.Range("I218").Value = 'sd'
$ws.Range("J218").Value = 'Statement-non-opinion'
$ws.Range("I219").Value = 'sv'
$ws.Range("J219").Value = 'Statement-opinion'
$ws.Range("I251").Value = 'ba'
$ws.Range("J251").Value = 'Appreciation'
$ws.Range("I253").Value = 'sv'
$ws.Range("J253").Value = 'Statement-opinion'
$ws.Range("I260").Value = 'aa'
$ws.Range("J260").Value = 'Agree/Accept'
$ws.Range("I263").Value = 'b'
$ws.Range("J263").Value = 'Acknowledge (Backchannel)'
$ws.Range("I265").Value = 'b'
$ws.Range("J265").Value = 'Acknowledge (Backchannel)'
$ws.Range("I268").Value = 'b'
$ws.Range("J268").Value = 'Acknowledge (Backchannel)'
$ws.Range("I271").Value = 'b'
$ws.Range("J271").Value = 'Acknowledge (Backchannel)'
$ws.Range("I281").Value = 'aa'
$ws.Range("J281").Value = 'Agree/Accept'
$ws.Range("I284").Value = 'ba'
$ws.Range("J284").Value = 'Appreciation'
$ws.Range("I297").Value = 'sd'
$ws.Range("J297").Value = 'Statement-non-opinion'
$ws.Range("I339").Value = 'aa'
$ws.Range("J339").Value = 'Agree/Accept'
$ws.Range("I342").Value = 'sv'
$ws.Range("J342").Value = 'Statement-opinion'
$ws.Range("I349").Value = 'aa'
$ws.Range("J349").Value = 'Agree/Accept'
$ws.Range("I355").Value = 'sd'
$ws.Range("J355").Value = 'Statement-non-opinion'
$ws.Range("I357").Value = 'sd'
$ws.Range("J357").Value = 'Statement-non-opinion'
$ws.Range("I358").Value = 'sd'
$ws.Range("J358").Value = 'Statement-non-opinion'
$ws.Range("I406").Value = 'sv'
$ws.Range("J406").Value = 'Statement-opinion'
$ws.Range("I416").Value = 'aa'
$ws.Range("J416").Value = 'Agree/Accept'
$ws.Range("I419").Value = 'sv'
$ws.Range("J419").Value = 'Statement-opinion'
$ws.Range("I420").Value = 'sv'
$ws.Range("J420").Value = 'Statement-opinion'
$ws.Range("I431").Value = 'sv'
$ws.Range("J431").Value = 'Statement-opinion'
$ws.Range("I435").Value = 'qy'
$ws.Range("J435").Value = 'Yes-No-Question'
$ws.Range("I445").Value = 'sd'
$ws.Range("J445").Value = 'Statement-non-opinion'
$ws.Range("I446").Value = 'sv'
$ws.Range("J446").Value = 'Statement-opinion'
$ws.Range("I472").Value = 'sd'
$ws.Range("J472").Value = 'Statement-non-opinion'
$ws.Range("I480").Value = 'b'
$ws.Range("J480").Value = 'Acknowledge (Backchannel)'
$ws.Range("I484").Value = '%'
$ws.Range("J484").Value = 'Uninterpretable'
$ws.Range("I485").Value = 'sd'
$ws.Range("J485").Value = 'Statement-non-opinion'
$ws.Range("I494").Value = 'sv'
$ws.Range("J494").Value = 'Statement-opinion'
$ws.Range("I498").Value = 'aa'
$ws.Range("J498").Value = 'Agree/Accept'
$ws.Range("I506").Value = '%'
$ws.Range("J506").Value = 'Uninterpretable'
$ws.Range("I511").Value = 'sd'
$ws.Range("J511").Value = 'Statement-non-opinion'
$ws.Range("I525").Value = 'aa'
$ws.Range("J525").Value = 'Agree/Accept'
$ws.Range("I526").Value = 'sd'
$ws.Range("J526").Value = 'Statement-non-opinion'
$ws.Range("I527").Value = 'aa'
$ws.Range("J527").Value = 'Agree/Accept'
